# Apply FHIR IG terminology / profile corrections to the Metadata sheet:
#  - Row 7  (Experimental) -> set Value column to "false"
#  - Row 8  (Date)         -> update the ISO timestamp
#  - Row 17 (Description)  -> set Value column to the CodeSystem description
#
# Note: a literal "false"/"true" string assigned through Range.Value is
# auto-coerced to a native Boolean by the engine (just like typing it into
# a cell in real Excel). Since we need it to remain a *text* value, we
# stage it on a scratch cell using a leading apostrophe (Excel's
# "force text" escape), copy that cell, and paste-special only the value
# into the destination so the destination keeps its original cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Row 7: Experimental = "false" (kept as text, not Boolean) ---
$scratch = $ws.Range("Z1000")
$scratch.Value = "'false"
$scratch.Copy()
$ws.Range("B7").PasteSpecial(-4163)   # xlPasteValues
$scratch.ClearContents()
$excel.CutCopyMode = 0

# --- Row 8: Date ---
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"

# --- Row 17: Description ---
$ws.Range("B17").Value = "Specific protocols and tests for VO2max estimation"

Write-Host "done"
